$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert 3 new rows at 8..10 (gender breakdown block) and a new column H
#    (2022 data). Doing the row insert first then the column insert mirrors
#    how the sheet grew (rows 1-7 gain an H cell that inherits the G column
#    style; rows 8-10 are brand new).
# ---------------------------------------------------------------------------
$ws.Rows("8:10").Insert()
$ws.Columns("H").Insert()

# ---------------------------------------------------------------------------
# 2. Fill in the new "2022" column header + data for the existing rows.
#    Cloning the format from the neighbouring "G" cell (PasteSpecial formats)
#    keeps every cell on the same cellXf as the rest of its row.
# ---------------------------------------------------------------------------
function Copy-Format($fromRef, $toRef) {
    $ws.Range($fromRef).Copy() | Out-Null
    $ws.Range($toRef).PasteSpecial(-4122) | Out-Null
}

Copy-Format "G3" "H3"
$ws.Range("H3").Value = 2022

Copy-Format "G4" "H4"
$ws.Range("H4").Value = 21.549331200908018

Copy-Format "G5" "H5"

Copy-Format "G6" "H6"
$ws.Range("H6").Value = 52.326989300763088

Copy-Format "G7" "H7"
$ws.Range("H7").Value = 3.6916333239218613

Write-Output "step1 ok"
